$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.335.33"
$ws.Range("E2").Value = "  +0.34%  "
$ws.Range("D3").Value = "1.669.05"
$ws.Range("E3").Value = "  +0.66%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.007"
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "220.73"
$ws.Range("E5").Value = "  +1.07%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5317"
$ws.Range("E6").Value = "  -0.08%  "
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2651"
$ws.Range("E8").Value = "  +0.84%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06368"
$ws.Range("E9").Value = "  +0.23%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.95"
$ws.Range("E10").Value = "  +2.26%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07844"
$ws.Range("E11").Value = "  +0.07%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.534"
$ws.Range("D13").Value = "1.677.85"
$ws.Range("E13").Value = "  +1.39%  "
$ws.Range("D14").Value = "1.897.89"
$ws.Range("E14").Value = "  +0.69%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5618"
$ws.Range("E15").Value = "  +1.85%  "
$ws.Range("D16").Value = "0.0₅8154"
$ws.Range("E16").Value = "  -0.45%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "65.90"
$ws.Range("E17").Value = "  +0.55%  "
$ws.Range("D18").Value = "26.329.30"
$ws.Range("E18").Value = "  +0.43%  "
$ws.Range("E19").Value = "  -0.10%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.726"
$ws.Range("E20").Value = "  +2.22%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "198.55"
$ws.Range("E21").Value = "  +3.41%  "
$ws.Range("E22").Value = "  +1.60%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.061"
$ws.Range("E23").Value = "  +0.48%  "
$ws.Range("E24").Value = "  -0.11%  "
$ws.Range("E25").Value = "  +1.98%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1218"
$ws.Range("E26").Value = "  -0.18%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.255"
$ws.Range("E27").Value = "  +0.33%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "16.15"
$ws.Range("E28").Value = "  +0.74%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.509"
$ws.Range("E30").Value = "  +2.02%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.285"
$ws.Range("E31").Value = "  +0.67%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.556"
$ws.Range("E32").Value = "  -0.28%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.326"
$ws.Range("E33").Value = "  +1.31%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.606"
$ws.Range("E34").Value = "  +0.24%  "
$ws.Range("B35").Value = "ARBITRUM"
$ws.Range("C35").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9632"
$ws.Range("E35").Value = "  +0.88%  "
$ws.Range("B36").Value = "MXToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.832"
$ws.Range("E36").Value = "  +0.53%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.432"
$ws.Range("E37").Value = "  +0.20%  "
$ws.Range("E38").Value = "  +0.61%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01618"
$ws.Range("E39").Value = "  +0.85%  "
$ws.Range("E40").Value = "  +2.39%  "
$ws.Range("E41").Value = "  +3.11%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8584"
$ws.Range("E42").Value = "  +0.69%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.007"
$ws.Range("E43").Value = "  -0.08%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "102.86"
$ws.Range("E44").Value = "  -1.79%  "
$ws.Range("D45").Value = "1.807.88"
$ws.Range("E45").Value = "  +0.54%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "58.62"
$ws.Range("E46").Value = "  +2.99%  "
$ws.Range("E47").Value = "  +0.59%  "
$ws.Range("B48").Value = "BabyDogeCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D48").Value = "0.0₈104"
$ws.Range("E48").Value = "  -0.85%  "
$ws.Range("B49").Value = "Mantle"
$ws.Range("C49").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.4403"
$ws.Range("E49").Value = "  +0.77%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.052"
$ws.Range("E50").Value = "  +1.74%  "
$ws.Range("E51").Value = "  -0.14%  "
